$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.598.04'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.228.43'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.12%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.47'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -8.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '296.82'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +11.00%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.88%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.29%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.43'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -8.69%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0917'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.31'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.20%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +10.28%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.62%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.561.86'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.238.98'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.499.82'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.81%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.32'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.06%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.92'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.53'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +21.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '229.85'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.20'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.05%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.65%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.46'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -11.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.24'

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.91%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.34'

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.26%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0897'

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.49%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.15'

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +11.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.30'

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.86%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0377'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.30%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.59%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.35'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.235'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.29%  '

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'Celestia'

$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.78'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -7.15%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'FirstDigitalUSD'

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.43'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -8.70%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.68%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.58'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.93%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.13'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.65'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +7.05%  '
